$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.971.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.100.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  -1.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09340"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.096.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.744"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.156"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06676"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.219"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.055.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.331"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.553"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.40"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1058"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.640"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.232"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.946"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.210"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02562"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06778"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.35%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.308"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6644"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.279"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.631"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000351"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.20%  "

$ws.Range("E49").Value = "  -3.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.54%  "
